$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top Gainers")

$ws.Cells.Item(2, 3).Value = 14.7054
$ws.Cells.Item(2, 4).Value = 10.6959
$ws.Cells.Item(2, 5).Value = 6.2133

$ws.Cells.Item(3, 3).Value = 14.5282
$ws.Cells.Item(3, 4).Value = 13.952
$ws.Cells.Item(3, 5).Value = 11.2046

$ws.Cells.Item(4, 3).Value = 11.8701
$ws.Cells.Item(4, 4).Value = 8.951599999999999
$ws.Cells.Item(4, 5).Value = 9.4025

$ws.Cells.Item(5, 3).Value = 11.625
$ws.Cells.Item(5, 4).Value = 16.1094
$ws.Cells.Item(5, 5).Value = 16.8684

$ws.Cells.Item(7, 3).Value = 9.325100000000001
$ws.Cells.Item(7, 4).Value = 7.6073
$ws.Cells.Item(7, 5).Value = 11.8999

$ws.Cells.Item(8, 3).Value = 7.9375
$ws.Cells.Item(8, 4).Value = 11.9352
$ws.Cells.Item(8, 5).Value = 14.3296

$ws.Cells.Item(9, 2).Value = 'BUTTERFLY'
$ws.Cells.Item(9, 3).Value = 7.5721
$ws.Cells.Item(9, 4).Value = 10.4809
$ws.Cells.Item(9, 5).Value = 13.0787

$ws.Cells.Item(10, 2).Value = 'FIVESTAR'
$ws.Cells.Item(10, 3).Value = 7.4481
$ws.Cells.Item(10, 4).Value = 7.4581
$ws.Cells.Item(10, 5).Value = 7.5382

$ws.Cells.Item(11, 2).Value = 'RPOWER'
$ws.Cells.Item(11, 3).Value = 7.326
$ws.Cells.Item(11, 4).Value = 3.7168
$ws.Cells.Item(11, 5).Value = 5.6808

$ws.Cells.Item(12, 2).Value = 'WALCHANNAG'
$ws.Cells.Item(12, 3).Value = 7.1817
$ws.Cells.Item(12, 4).Value = 4.7064
$ws.Cells.Item(12, 5).Value = -4.5406

$ws.Cells.Item(13, 3).Value = 7.1639
$ws.Cells.Item(13, 4).Value = 7.941
$ws.Cells.Item(13, 5).Value = 7.4938

$ws.Cells.Item(14, 3).Value = 6.6355
$ws.Cells.Item(14, 4).Value = 7.2614
$ws.Cells.Item(14, 5).Value = 16.1642

$ws.Cells.Item(16, 3).Value = 6.289
$ws.Cells.Item(16, 4).Value = 11.2988
$ws.Cells.Item(16, 5).Value = 4.4505

$ws.Cells.Item(17, 2).Value = 'SANDUMA'
$ws.Cells.Item(17, 3).Value = 6.0146
$ws.Cells.Item(17, 4).Value = 3.5288
$ws.Cells.Item(17, 5).Value = 32.0521

$ws.Cells.Item(18, 2).Value = 'UTKARSHBNK'
$ws.Cells.Item(18, 3).Value = 5.9548
$ws.Cells.Item(18, 4).Value = -4.9286
$ws.Cells.Item(18, 5).Value = -1.6206

$ws.Cells.Item(19, 2).Value = 'POKARNA'
$ws.Cells.Item(19, 3).Value = 5.7265
$ws.Cells.Item(19, 4).Value = -1.687
$ws.Cells.Item(19, 5).Value = 18.4942

$ws.Cells.Item(20, 2).Value = 'JISLJALEQS'
$ws.Cells.Item(20, 3).Value = 5.673
$ws.Cells.Item(20, 4).Value = 4.9687
$ws.Cells.Item(20, 5).Value = -1.1037

$ws.Cells.Item(21, 2).Value = 'GENUSPOWER'
$ws.Cells.Item(21, 3).Value = 5.5237
$ws.Cells.Item(21, 4).Value = 3.8225
$ws.Cells.Item(21, 5).Value = 0.7418

$ws.Cells.Item(22, 2).Value = 'CELLO'
$ws.Cells.Item(22, 3).Value = 5.499
$ws.Cells.Item(22, 4).Value = 4.3364
$ws.Cells.Item(22, 5).Value = 14.2448

$ws.Cells.Item(23, 2).Value = 'ADANIENSOL'
$ws.Cells.Item(23, 3).Value = 5.4711
$ws.Cells.Item(23, 4).Value = 2.8965
$ws.Cells.Item(23, 5).Value = 11.3901

$ws.Cells.Item(24, 2).Value = 'VAIBHAVGBL'
$ws.Cells.Item(24, 3).Value = 5.456
$ws.Cells.Item(24, 4).Value = 5.9817
$ws.Cells.Item(24, 5).Value = 12.3041

$ws.Cells.Item(25, 2).Value = 'EPACKPEB'
$ws.Cells.Item(25, 3).Value = 5.326
$ws.Cells.Item(25, 4).Value = -1.184
$ws.Cells.Item(25, 5).Value = 'N/A'

$ws.Cells.Item(26, 2).Value = 'ABDL'
$ws.Cells.Item(26, 3).Value = 5.2786
$ws.Cells.Item(26, 4).Value = 4.1805
$ws.Cells.Item(26, 5).Value = 26.886

$ws.Cells.Item(27, 3).Value = 5.1903
$ws.Cells.Item(27, 4).Value = 8.093400000000001
$ws.Cells.Item(27, 5).Value = 8.511900000000001

$ws.Cells.Item(28, 2).Value = 'GRAPHITE'
$ws.Cells.Item(28, 3).Value = 5.0568
$ws.Cells.Item(28, 4).Value = 11.4476
$ws.Cells.Item(28, 5).Value = 11.6482

$ws.Cells.Item(29, 2).Value = 'MEGASOFT'
$ws.Cells.Item(29, 3).Value = 4.9974
$ws.Cells.Item(29, 4).Value = 15.7588
$ws.Cells.Item(29, 5).Value = 33.5271

$ws.Cells.Item(30, 2).Value = 'PROZONER'
$ws.Cells.Item(30, 3).Value = 4.9921
$ws.Cells.Item(30, 4).Value = 15.7468
$ws.Cells.Item(30, 5).Value = 36.095

$ws.Cells.Item(31, 2).Value = 'STALLION'
$ws.Cells.Item(31, 3).Value = 4.9914
$ws.Cells.Item(31, 4).Value = -5.2229
$ws.Cells.Item(31, 5).Value = 21.4391

$ws.Cells.Item(32, 2).Value = 'INDOTHAI'
$ws.Cells.Item(32, 3).Value = 4.9883
$ws.Cells.Item(32, 4).Value = 4.7163
$ws.Cells.Item(32, 5).Value = 43.9974

$ws.Cells.Item(33, 2).Value = 'ATGL'
$ws.Cells.Item(33, 3).Value = 4.8712
$ws.Cells.Item(33, 4).Value = 4.6101
$ws.Cells.Item(33, 5).Value = 4.15

$ws.Cells.Item(35, 2).Value = 'SURYAROSNI'
$ws.Cells.Item(35, 3).Value = 4.7517
$ws.Cells.Item(35, 4).Value = 11.1405
$ws.Cells.Item(35, 5).Value = 2.7943

$ws.Cells.Item(36, 2).Value = 'DATAMATICS'
$ws.Cells.Item(36, 3).Value = 4.722
$ws.Cells.Item(36, 4).Value = 7.1326
$ws.Cells.Item(36, 5).Value = 15.5329

$ws.Cells.Item(37, 3).Value = 4.6754
$ws.Cells.Item(37, 4).Value = 3.4784
$ws.Cells.Item(37, 5).Value = 2.2843

$ws.Cells.Item(38, 3).Value = 4.504
$ws.Cells.Item(38, 4).Value = 8.2417
$ws.Cells.Item(38, 5).Value = 15.4758

$ws.Cells.Item(39, 3).Value = 4.4554
$ws.Cells.Item(39, 4).Value = 6.9204
$ws.Cells.Item(39, 5).Value = -3.074

$ws.Cells.Item(40, 2).Value = 'GMBREW'
$ws.Cells.Item(40, 3).Value = 4.4369
$ws.Cells.Item(40, 4).Value = -0.0158
$ws.Cells.Item(40, 5).Value = 79.95440000000001

$ws.Cells.Item(41, 2).Value = 'CMSINFO'
$ws.Cells.Item(41, 3).Value = 4.4237
$ws.Cells.Item(41, 4).Value = 3.1952
$ws.Cells.Item(41, 5).Value = 3.4025

$ws.Cells.Item(42, 2).Value = 'SAMBHV'
$ws.Cells.Item(42, 3).Value = 4.4166
$ws.Cells.Item(42, 4).Value = 2.9017
$ws.Cells.Item(42, 5).Value = 5.4515

$ws.Cells.Item(43, 2).Value = 'PDSL'
$ws.Cells.Item(43, 3).Value = 4.3814
$ws.Cells.Item(43, 4).Value = 2.4096
$ws.Cells.Item(43, 5).Value = 8.2037

$ws.Cells.Item(44, 2).Value = 'BAJAJINDEF'
$ws.Cells.Item(44, 3).Value = 4.265
$ws.Cells.Item(44, 4).Value = 3.1883
$ws.Cells.Item(44, 5).Value = 10.1861

$ws.Cells.Item(45, 2).Value = 'SUNFLAG'
$ws.Cells.Item(45, 3).Value = 4.1675
$ws.Cells.Item(45, 4).Value = 4.504
$ws.Cells.Item(45, 5).Value = 4.8027

$ws.Cells.Item(46, 2).Value = 'STLTECH'
$ws.Cells.Item(46, 3).Value = 4.1667
$ws.Cells.Item(46, 4).Value = 1.1741
$ws.Cells.Item(46, 5).Value = 7.2658

$ws.Cells.Item(47, 2).Value = 'PROSTARM'
$ws.Cells.Item(47, 3).Value = 4.1532
$ws.Cells.Item(47, 4).Value = 1.2644
$ws.Cells.Item(47, 5).Value = -7.6891

$ws.Cells.Item(48, 2).Value = 'SGMART'
$ws.Cells.Item(48, 3).Value = 4.1185
$ws.Cells.Item(48, 4).Value = 8.097799999999999
$ws.Cells.Item(48, 5).Value = 2.3856

$ws.Cells.Item(49, 2).Value = 'GPIL'
$ws.Cells.Item(49, 3).Value = 4.0599
$ws.Cells.Item(49, 4).Value = 6.2282
$ws.Cells.Item(49, 5).Value = 14.3342

$ws.Cells.Item(50, 2).Value = 'LLOYDSENT'
$ws.Cells.Item(50, 3).Value = 4.0223
$ws.Cells.Item(50, 4).Value = 1.3058
$ws.Cells.Item(50, 5).Value = 10.6571

$ws.Cells.Item(51, 2).Value = 'RHIM'
$ws.Cells.Item(51, 3).Value = 3.9338
$ws.Cells.Item(51, 4).Value = 3.5058
$ws.Cells.Item(51, 5).Value = 5.4661

$ws.Cells.Item(52, 2).Value = 'GPPL'
$ws.Cells.Item(52, 3).Value = 3.9059
$ws.Cells.Item(52, 4).Value = 2.9027
$ws.Cells.Item(52, 5).Value = 4.5371

$ws.Cells.Item(53, 3).Value = 3.8602
$ws.Cells.Item(53, 4).Value = 0.5569
$ws.Cells.Item(53, 5).Value = 32.7208

$ws.Cells.Item(54, 2).Value = 'ICRA'
$ws.Cells.Item(54, 3).Value = 3.8382
$ws.Cells.Item(54, 4).Value = 4.5193
$ws.Cells.Item(54, 5).Value = 2.9222

$ws.Cells.Item(55, 2).Value = 'RECLTD'
$ws.Cells.Item(55, 3).Value = 3.7674
$ws.Cells.Item(55, 4).Value = 2.7509
$ws.Cells.Item(55, 5).Value = 2.682

$ws.Cells.Item(56, 2).Value = 'TCI'
$ws.Cells.Item(56, 3).Value = 3.7609
$ws.Cells.Item(56, 4).Value = 3.6647
$ws.Cells.Item(56, 5).Value = 4.1649

$ws.Cells.Item(57, 2).Value = 'NBCC'
$ws.Cells.Item(57, 3).Value = 3.7259
$ws.Cells.Item(57, 4).Value = 2.4443
$ws.Cells.Item(57, 5).Value = 6.8547

$ws.Cells.Item(58, 2).Value = 'SRM'
$ws.Cells.Item(58, 3).Value = 3.7168
$ws.Cells.Item(58, 4).Value = 3.4086
$ws.Cells.Item(58, 5).Value = 4.3193

$ws.Cells.Item(59, 2).Value = 'MRPL'
$ws.Cells.Item(59, 3).Value = 3.6691
$ws.Cells.Item(59, 4).Value = 9.084
$ws.Cells.Item(59, 5).Value = 19.3689

$ws.Cells.Item(60, 2).Value = 'STAR'
$ws.Cells.Item(60, 3).Value = 3.6487
$ws.Cells.Item(60, 4).Value = 3.5787
$ws.Cells.Item(60, 5).Value = 2.815

$ws.Cells.Item(61, 2).Value = 'HCC'
$ws.Cells.Item(61, 3).Value = 3.6377
$ws.Cells.Item(61, 4).Value = 2.5406
$ws.Cells.Item(61, 5).Value = 7.2721

$ws.Cells.Item(62, 2).Value = 'BLACKBUCK'
$ws.Cells.Item(62, 3).Value = 3.6359
$ws.Cells.Item(62, 4).Value = 2.2526
$ws.Cells.Item(62, 5).Value = 8.2196

$ws.Cells.Item(63, 2).Value = 'VINCOFE'
$ws.Cells.Item(63, 3).Value = 3.6218
$ws.Cells.Item(63, 4).Value = 10.4782
$ws.Cells.Item(63, 5).Value = 8.8552

$ws.Cells.Item(64, 3).Value = 3.6086
$ws.Cells.Item(64, 4).Value = 2.1434
$ws.Cells.Item(64, 5).Value = 4.884

$ws.Cells.Item(65, 2).Value = 'ASHAPURMIN'
$ws.Cells.Item(65, 3).Value = 3.5558
$ws.Cells.Item(65, 4).Value = 6.206
$ws.Cells.Item(65, 5).Value = 1.9969

$ws.Cells.Item(66, 2).Value = 'MSTCLTD'
$ws.Cells.Item(66, 3).Value = 3.555
$ws.Cells.Item(66, 4).Value = 3.4871
$ws.Cells.Item(66, 5).Value = 15.869

$ws.Cells.Item(67, 3).Value = 3.5545
$ws.Cells.Item(67, 4).Value = 5.6417
$ws.Cells.Item(67, 5).Value = 14.093

$ws.Cells.Item(68, 2).Value = 'INDORAMA'
$ws.Cells.Item(68, 3).Value = 3.5365
$ws.Cells.Item(68, 4).Value = 2.7019
$ws.Cells.Item(68, 5).Value = 13.8319

$ws.Cells.Item(69, 2).Value = 'PRAKASH'
$ws.Cells.Item(69, 3).Value = 3.535
$ws.Cells.Item(69, 4).Value = 4.4385
$ws.Cells.Item(69, 5).Value = 1.1861

$ws.Cells.Item(70, 2).Value = 'MAITHANALL'
$ws.Cells.Item(70, 3).Value = 3.5328
$ws.Cells.Item(70, 4).Value = 2.5747
$ws.Cells.Item(70, 5).Value = 1.8468

$ws.Cells.Item(71, 2).Value = 'RAJRATAN'
$ws.Cells.Item(71, 3).Value = 3.5239
$ws.Cells.Item(71, 4).Value = 1.0712
$ws.Cells.Item(71, 5).Value = 27.12

$ws.Cells.Item(72, 2).Value = 'DCMSHRIRAM'
$ws.Cells.Item(72, 3).Value = 3.5092
$ws.Cells.Item(72, 4).Value = 10.0678
$ws.Cells.Item(72, 5).Value = 17.4318

$ws.Cells.Item(75, 2).Value = 'PENIND'
$ws.Cells.Item(75, 3).Value = 3.4578
$ws.Cells.Item(75, 4).Value = 2.7376
$ws.Cells.Item(75, 5).Value = 12.0973

$ws.Cells.Item(76, 2).Value = 'DCW'
$ws.Cells.Item(76, 3).Value = 3.4203
$ws.Cells.Item(76, 4).Value = 1.9925
$ws.Cells.Item(76, 5).Value = -4.2844

Write-Host "Applied 2025-10-29 08:35 data update to Top Gainers sheet"
